$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for "Zapallo italiano" at Vega Monumental
# Concepción was inserted ahead of the former row 150, pushing every
# subsequent record down by one row (old 150..205 -> new 151..206).
$ws.Rows.Item(150).Insert()

$ws.Range("A150").Value = 11
$ws.Range("B150").Value = "Vega Monumental Concepción"
$ws.Range("C150").Value = "Bíobío"
$ws.Range("D150").Value = 45006
$ws.Range("E150").Value = 8
$ws.Range("F150").Value = 100112032
$ws.Range("G150").Value = "Zapallo italiano"
$ws.Range("H150").Value = "Sin especificar"
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 100
$ws.Range("K150").Value = 8500
$ws.Range("L150").Value = 9000
$ws.Range("M150").Value = 8750
$ws.Range("N150").Value = "$/caja 50 unidades"
$ws.Range("O150").Value = "Región de Arica y Parinacota"
$ws.Range("P150").Value = 175
$ws.Range("Q150").Value = 50
$ws.Range("R150").Value = "Hortaliza"
